$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 7 data (Both Positive and Negative Scenarios - data driven test) ---
# A7 stays empty but copies the formatting used by A3:A6 (style idx 3)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

# B7 gets the "test324" value, formatted like B6 (hyperlink-like style idx 5)
$ws.Range("B6").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B7").Value = "test324"

# C7 gets the "Invalid" value (reuses existing shared string), formatted like C6 (style idx 6)
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C7").Value = "Invalid"

# Match row height used by the rest of the data rows (21)
$ws.Rows(7).RowHeight = 21

# --- Update the active selection shown in the sheet view ---
$ws.Range("A10").Select()

# --- Bump the cached max row outline level from 5 to 6 without leaving any row ---
# with an actual outlineLevel attribute (mirrors the stray metadata-only change in the diff)
$ws.Rows("50:50").OutlineLevel = 6
$ws.Rows("50:50").Delete()

Write-Host "Edit applied"
